$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("CY2").Value = "P"
